$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '44.929.70'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -5.88%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.675.88'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.05%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.11%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.39'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.46%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '97.52'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.88%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.592'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.62%  '

# Row 8
$ws.Range('E8').Value = '  +0.02%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.574'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.26%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '37.86'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.96%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0838'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.59%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.00'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.53%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.082.08'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.61%  '

# Row 14
$ws.Range('E14').Value = '  +0.20%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.686.73'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.75%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.922'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.43%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '15.00'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.77%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '44.884.61'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -6.00%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.83'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.09%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000100'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.18%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.62'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.59%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '75.05'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.69%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '277.47'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.15%  '

# Row 24
$ws.Range('B24').Value = 'ImmutableX'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.28'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.81%  '

# Row 25
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.01'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.30%  '

# Row 26
$ws.Range('E26').Value = '  +1.15%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.998'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.34%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.51'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.67%  '

# Row 29
$ws.Range('E29').Value = '  -3.80%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '37.57'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.68%  '

# Row 31
$ws.Range('B31').Value = 'LidoDAOToken'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.78'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.81%  '

# Row 32
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.16'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.63%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.34'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.44%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '153.42'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.36%  '

# Row 35
$ws.Range('B35').Value = 'WEMIXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.81'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.07%  '

# Row 36
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0832'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.81%  '

# Row 37
$ws.Range('E37').Value = '  -7.27%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '25.02'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +11.09%  '

# Row 39
$ws.Range('E39').Value = '  -0.46%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '15.94'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.51%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.58'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.22%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0322'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.08%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.148.14'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.48%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.92'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -8.33%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.996'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.23%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '92.28'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.23%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.42'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.62%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.940.14'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.06%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '110.28'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.11%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.62'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.79%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.197'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.86%  '
